$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 150 (existing rows 150..252 shift down to 152..254).
$ws.Range("A150:T151").EntireRow.Insert()

# --- New row 150 ---
$ws.Cells.Item(150, 1).Value = 7
$ws.Cells.Item(150, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(150, 3).Value = "Ñuble"
$ws.Cells.Item(150, 4).Value = 45062
$ws.Cells.Item(150, 5).Value = 16
$ws.Cells.Item(150, 6).Value = "Fruta"
$ws.Cells.Item(150, 7).Value = 100101
$ws.Cells.Item(150, 8).Value = "Berries"
$ws.Cells.Item(150, 9).Value = 100101007
$ws.Cells.Item(150, 10).Value = "Kiwi"
$ws.Cells.Item(150, 11).Value = "Hayward"
$ws.Cells.Item(150, 12).Value = "Especial"
$ws.Cells.Item(150, 13).Value = 50
$ws.Cells.Item(150, 14).Value = 12000
$ws.Cells.Item(150, 15).Value = 12000
$ws.Cells.Item(150, 16).Value = 12000
$ws.Cells.Item(150, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(150, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(150, 19).Value = 667
$ws.Cells.Item(150, 20).Value = 18

# --- New row 151 ---
$ws.Cells.Item(151, 1).Value = 7
$ws.Cells.Item(151, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(151, 3).Value = "Ñuble"
$ws.Cells.Item(151, 4).Value = 45062
$ws.Cells.Item(151, 5).Value = 16
$ws.Cells.Item(151, 6).Value = "Fruta"
$ws.Cells.Item(151, 7).Value = 100101
$ws.Cells.Item(151, 8).Value = "Berries"
$ws.Cells.Item(151, 9).Value = 100101007
$ws.Cells.Item(151, 10).Value = "Kiwi"
$ws.Cells.Item(151, 11).Value = "Hayward"
$ws.Cells.Item(151, 12).Value = "Primera"
$ws.Cells.Item(151, 13).Value = 40
$ws.Cells.Item(151, 14).Value = 10000
$ws.Cells.Item(151, 15).Value = 10000
$ws.Cells.Item(151, 16).Value = 10000
$ws.Cells.Item(151, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(151, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(151, 19).Value = 556
$ws.Cells.Item(151, 20).Value = 18
